$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = 0.0767925000000001
$ws1.Range("E2").Value = 0.4072877638625012
$ws1.Range("G2").Value = 0.2494892361374989
$ws1.Range("I2").Value = 0.5224279530531742
$ws1.Range("L2").Value = 0.5729273469468261
$ws1.Range("N2").Value = 11.81230257363055
$ws1.Range("O2").Value = 2.569167174609038

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 0.09263141357332233
$ws2.Range("E2").Value = 0.2374014999999999
$ws2.Range("I2").Value = 0.4472572469468257
$ws2.Range("L2").Value = 0.1550176530531738
$ws2.Range("N2").Value = 5.761745064355102
$ws2.Range("O2").Value = 1.463791048667069

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 0.0509340911268795
$ws3.Range("B2").Value = 0.05174907978289767
$ws3.Range("E2").Value = 0.2384590155169009
$ws3.Range("I2").Value = 0.4896953000000001
$ws3.Range("L2").Value = 0
$ws3.Range("M2").Value = 0
$ws3.Range("N2").Value = 8.485441261999803
$ws3.Range("O2").Value = 5.722779367126569
